# Updated cryptos list on Sun Jul  7 04:45:22 UTC 2024 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) for each coin row.
# D-column values are prefixed with a leading apostrophe so Excel stores
# them as text (preserving formats like "64.80" / "57.769.69") instead of
# auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''57.769.69'
$ws.Range('E2').Value = '  +2.91%  '
$ws.Range('D3').Value = '''3.036.54'
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''512.51'
$ws.Range('E5').Value = '  +2.47%  '
$ws.Range('D6').Value = '''139.96'
$ws.Range('E6').Value = '  +4.49%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +3.37%  '
$ws.Range('E9').Value = '  +2.33%  '
$ws.Range('E10').Value = '  +4.45%  '
$ws.Range('E11').Value = '  +5.30%  '
$ws.Range('D12').Value = '''3.564.39'
$ws.Range('E12').Value = '  +2.51%  '
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range('D14').Value = '''26.65'
$ws.Range('E14').Value = '  +6.13%  '
$ws.Range('E15').Value = '  +11.52%  '
$ws.Range('D16').Value = '''57.777.94'
$ws.Range('E16').Value = '  +2.92%  '
$ws.Range('E17').Value = '  +9.90%  '
$ws.Range('D18').Value = '''3.038.28'
$ws.Range('E18').Value = '  +2.50%  '
$ws.Range('D19').Value = '''12.96'
$ws.Range('E19').Value = '  +5.74%  '
$ws.Range('E20').Value = '  +4.05%  '
$ws.Range('D21').Value = '''335.34'
$ws.Range('E21').Value = '  +4.43%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  +6.75%  '
$ws.Range('D25').Value = '''64.80'
$ws.Range('E25').Value = '  +4.95%  '
$ws.Range('D26').Value = '''0.167'
$ws.Range('E26').Value = '  +3.92%  '
$ws.Range('D27').Value = '''0.999'
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('D28').Value = '''0.0₃0935'
$ws.Range('E28').Value = '  +6.37%  '
$ws.Range('D29').Value = '''6.81'
$ws.Range('E29').Value = '  +6.01%  '
$ws.Range('D30').Value = '''7.44'
$ws.Range('E30').Value = '  +10.52%  '
$ws.Range('E31').Value = '  +4.38%  '
$ws.Range('E32').Value = '  +3.18%  '
$ws.Range('E33').Value = '  +2.84%  '
$ws.Range('D34').Value = '''156.13'
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('E35').Value = '  +6.81%  '
$ws.Range('E36').Value = '  +6.45%  '
$ws.Range('E37').Value = '  +2.25%  '
$ws.Range('D38').Value = '''24.73'
$ws.Range('E38').Value = '  +8.09%  '
$ws.Range('D39').Value = '''0.0686'
$ws.Range('E39').Value = '  +2.59%  '
$ws.Range('D40').Value = '''3.073.44'
$ws.Range('E40').Value = '  +2.38%  '
$ws.Range('D41').Value = '''37.52'
$ws.Range('E41').Value = '  +3.82%  '
$ws.Range('D42').Value = '''3.88'
$ws.Range('E42').Value = '  +10.01%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('E44').Value = '  +3.49%  '
$ws.Range('D45').Value = '''2.303.05'
$ws.Range('E45').Value = '  +2.86%  '
$ws.Range('E46').Value = '  +3.31%  '
$ws.Range('E47').Value = '  +2.02%  '
$ws.Range('D48').Value = '''6.03'
$ws.Range('E48').Value = '  +5.31%  '
$ws.Range('E49').Value = '  +3.42%  '
$ws.Range('D50').Value = '''19.55'
$ws.Range('E50').Value = '  +3.95%  '
$ws.Range('D51').Value = '''1.84'
$ws.Range('E51').Value = '  -3.58%  '
